$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data (D: Price, E: Volume(1h))
# Values are forced to remain text (matching the source inlineStr cells) by
# temporarily setting a text NumberFormat, then clearing the format afterward
# so no stray style index is left behind on the cell.
$updates = @(
    @{ Cell = 'D2'; Value = '63.402.12' }
    @{ Cell = 'E2'; Value = '  -3.65%  ' }
    @{ Cell = 'D3'; Value = '2.591.13' }
    @{ Cell = 'E3'; Value = '  -2.75%  ' }
    @{ Cell = 'E4'; Value = '  -0.01%  ' }
    @{ Cell = 'D5'; Value = '571.75' }
    @{ Cell = 'E5'; Value = '  -4.50%  ' }
    @{ Cell = 'D6'; Value = '155.09' }
    @{ Cell = 'E6'; Value = '  -1.93%  ' }
    @{ Cell = 'D7'; Value = '1.00' }
    @{ Cell = 'E7'; Value = '  +0.05%  ' }
    @{ Cell = 'D8'; Value = '0.622' }
    @{ Cell = 'E8'; Value = '  -5.09%  ' }
    @{ Cell = 'E9'; Value = '  -6.94%  ' }
    @{ Cell = 'E10'; Value = '  -0.13%  ' }
    @{ Cell = 'D12'; Value = '0.157' }
    @{ Cell = 'E12'; Value = '  -0.34%  ' }
    @{ Cell = 'E13'; Value = '  -2.81%  ' }
    @{ Cell = 'D14'; Value = '3.056.43' }
    @{ Cell = 'E14'; Value = '  -2.68%  ' }
    @{ Cell = 'E15'; Value = '  -7.99%  ' }
    @{ Cell = 'D16'; Value = '63.239.04' }
    @{ Cell = 'E16'; Value = '  -3.71%  ' }
    @{ Cell = 'D17'; Value = '2.605.14' }
    @{ Cell = 'E17'; Value = '  -2.25%  ' }
    @{ Cell = 'D18'; Value = '11.98' }
    @{ Cell = 'E18'; Value = '  -4.83%  ' }
    @{ Cell = 'D19'; Value = '7.53' }
    @{ Cell = 'E19'; Value = '  +0.64%  ' }
    @{ Cell = 'E20'; Value = '  -5.68%  ' }
    @{ Cell = 'D21'; Value = '342.20' }
    @{ Cell = 'E21'; Value = '  -2.66%  ' }
    @{ Cell = 'E22'; Value = '  +0.06%  ' }
    @{ Cell = 'D23'; Value = '67.18' }
    @{ Cell = 'E23'; Value = '  -3.81%  ' }
    @{ Cell = 'E24'; Value = '  -0.43%  ' }
    @{ Cell = 'E25'; Value = '  -3.63%  ' }
    @{ Cell = 'D26'; Value = '9.11' }
    @{ Cell = 'E26'; Value = '  -5.88%  ' }
    @{ Cell = 'D27'; Value = '578.96' }
    @{ Cell = 'E27'; Value = '  +1.30%  ' }
    @{ Cell = 'E28'; Value = '  -3.70%  ' }
    @{ Cell = 'D29'; Value = '0.999' }
    @{ Cell = 'E29'; Value = '  -0.17%  ' }
    @{ Cell = 'E30'; Value = '  -1.57%  ' }
    @{ Cell = 'D31'; Value = '7.87' }
    @{ Cell = 'E31'; Value = '  -4.08%  ' }
    @{ Cell = 'E33'; Value = '  -5.80%  ' }
    @{ Cell = 'D34'; Value = '6.55' }
    @{ Cell = 'E34'; Value = '  -2.26%  ' }
    @{ Cell = 'D35'; Value = '5.45' }
    @{ Cell = 'E35'; Value = '  -2.47%  ' }
    @{ Cell = 'E36'; Value = '  -4.90%  ' }
    @{ Cell = 'E37'; Value = '  +0.05%  ' }
    @{ Cell = 'D38'; Value = '19.70' }
    @{ Cell = 'E38'; Value = '  -4.52%  ' }
    @{ Cell = 'D39'; Value = '154.60' }
    @{ Cell = 'E39'; Value = '  -0.39%  ' }
    @{ Cell = 'E40'; Value = '  -4.65%  ' }
    @{ Cell = 'E41'; Value = '  -0.01%  ' }
    @{ Cell = 'D42'; Value = '41.24' }
    @{ Cell = 'E42'; Value = '  -3.57%  ' }
    @{ Cell = 'E43'; Value = '  +6.66%  ' }
    @{ Cell = 'D44'; Value = '155.80' }
    @{ Cell = 'E44'; Value = '  -3.48%  ' }
    @{ Cell = 'E45'; Value = '  -4.49%  ' }
    @{ Cell = 'D46'; Value = '23.08' }
    @{ Cell = 'E46'; Value = '  +0.20%  ' }
    @{ Cell = 'E47'; Value = '  -5.50%  ' }
    @{ Cell = 'E48'; Value = '  -3.17%  ' }
    @{ Cell = 'E49'; Value = '  -1.80%  ' }
    @{ Cell = 'D50'; Value = '0.0245' }
    @{ Cell = 'E50'; Value = '  -5.33%  ' }
    @{ Cell = 'D51'; Value = '18.75' }
    @{ Cell = 'E51'; Value = '  -5.34%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.ClearFormats()
}
